$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row for the new "OutputAnalysis2" table (AA1:AI1)
$headers = @("Name","Value","OurCount","TheirCount","OurSum","TheirSum","OurCount%","TheirCount%","Differential")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, 27 + $i).Value = $headers[$i]
}

# Data rows (AA2:AI7)
$ws.Range("AA2").Value = "carName"
$ws.Range("AB2").Value = "Fennec"
$ws.Range("AC2").Value = 507
$ws.Range("AD2").Value = 102
$ws.Range("AE2").Value = 912
$ws.Range("AF2").Value = 134
$ws.Range("AG2").Value = 0.555921052631579
$ws.Range("AH2").Value = 0.7611940298507462
$ws.Range("AI2").Value = -0.2052729772191673

$ws.Range("AA3").Value = "carName"
$ws.Range("AB3").Value = "Octane"
$ws.Range("AC3").Value = 405
$ws.Range("AD3").Value = 32
$ws.Range("AE3").Value = 912
$ws.Range("AF3").Value = 134
$ws.Range("AG3").Value = 0.4440789473684211
$ws.Range("AH3").Value = 0.2388059701492537
$ws.Range("AI3").Value = 0.2052729772191673

$ws.Range("AA4").Value = "mvp"
$ws.Range("AB4").Value = 0
$ws.Range("AC4").Value = 590
$ws.Range("AD4").Value = 85
$ws.Range("AE4").Value = 912
$ws.Range("AF4").Value = 134
$ws.Range("AG4").Value = 0.6469298245614035
$ws.Range("AH4").Value = 0.6343283582089553
$ws.Range("AI4").Value = 0.01260146635244819

$ws.Range("AA5").Value = "mvp"
$ws.Range("AB5").Value = 1
$ws.Range("AC5").Value = 322
$ws.Range("AD5").Value = 49
$ws.Range("AE5").Value = 912
$ws.Range("AF5").Value = 134
$ws.Range("AG5").Value = 0.3530701754385965
$ws.Range("AH5").Value = 0.3656716417910448
$ws.Range("AI5").Value = -0.0126014663524483

$ws.Range("AA6").Value = "scoredFirst"
$ws.Range("AB6").Value = $false
$ws.Range("AC6").Value = 477
$ws.Range("AD6").Value = 44
$ws.Range("AE6").Value = 677
$ws.Range("AF6").Value = 61
$ws.Range("AG6").Value = 0.7045790251107829
$ws.Range("AH6").Value = 0.7213114754098361
$ws.Range("AI6").Value = -0.01673245029905324

$ws.Range("AA7").Value = "scoredFirst"
$ws.Range("AB7").Value = $true
$ws.Range("AC7").Value = 200
$ws.Range("AD7").Value = 17
$ws.Range("AE7").Value = 677
$ws.Range("AF7").Value = 61
$ws.Range("AG7").Value = 0.2954209748892171
$ws.Range("AH7").Value = 0.2786885245901639
$ws.Range("AI7").Value = 0.01673245029905324

# Create the new table "OutputAnalysis2" over AA1:AI7, matching the
# header names already written into the sheet above.
$lo = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $ws.Range("AA1:AI7"), $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$lo.Name = "OutputAnalysis2"
$lo.TableStyle = "TableStyleMedium9"
$lo.ShowTableStyleColumnStripes = $true
